$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Bring over the formatting (styles/number formats) of column F into the
# new column G by copying the whole column range and pasting it in place.
$ws.Range("F1:F19").Copy()
$ws.Range("G1:G19").PasteSpecial(-4122)

# Re-set the actual values explicitly, since pasting a literal 0 can be
# dropped by copy/paste - setting Value directly guarantees the <v>0</v>
# is written out for every data/total row.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Header text for the new column
$ws.Range("G1").Value = "PRESUPUESTO"

# Column G width -> stored OOXML width of 17 (Excel's ColumnWidth property
# is offset from the stored width by the default cell padding, ~0.8333).
$ws.Columns.Item(7).ColumnWidth = 16.166666666666668

Write-Host "Column G (PRESUPUESTO) added to VENTA MENSUAL"
